$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 10359.6
$ws.Range("I86").Value = 7999.5
$ws.Range("K86").Value = 7999.5
$ws.Range("M86").Value = -6876.5
$ws.Range("H89").Value = 10359.6
$ws.Range("I89").Value = 7999.5
$ws.Range("K89").Value = 39997.5
$ws.Range("M89").Value = -34381.5
$ws.Range("H100").Value = 1754.4166
$ws.Range("I100").Value = 1639.7778
$ws.Range("J100").Value = 2098.3333
$ws.Range("K100").Value = 1639.7778
$ws.Range("L100").Value = 2098.3333
$ws.Range("M100").Value = -1098.7778
$ws.Range("N100").Value = -3180.3333
$ws.Range("H106").Value = 9000
$ws.Range("I106").Value = 9000
$ws.Range("K106").Value = 9000
$ws.Range("M106").Value = -8369
$ws.Range("H116").Value = 4144.364
$ws.Range("I116").Value = 3364.6667
$ws.Range("K116").Value = 3364.6667
$ws.Range("M116").Value = 77.33329999999978
$ws.Range("H129").Value = 1264566.5
$ws.Range("I129").Value = 1189.091
$ws.Range("K129").Value = 3567.273
$ws.Range("M129").Value = 1432.727
$ws.Range("H137").Value = 2170.4736
$ws.Range("I137").Value = 1970.1333
$ws.Range("J137").Value = 2921.75
$ws.Range("K137").Value = 5910.3999
$ws.Range("L137").Value = 8765.25
$ws.Range("M137").Value = -3360.3999
$ws.Range("N137").Value = -13865.25
$ws.Range("H138").Value = 3989.8076
$ws.Range("J138").Value = 4407
$ws.Range("L138").Value = 13221
$ws.Range("N138").Value = -23501

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6999.9
$ws.Range("I2").Value = 6249.875
$ws.Range("K2").Value = 6249.875
$ws.Range("M2").Value = -6136.875
$ws.Range("H32").Value = 3107.6233
$ws.Range("I32").Value = 2951.2769
$ws.Range("J32").Value = 5648.25
$ws.Range("K32").Value = 2951.2769
$ws.Range("L32").Value = 5648.25
$ws.Range("M32").Value = -2664.2769
$ws.Range("N32").Value = -6222.25
$ws.Range("H106").Value = 21106
$ws.Range("J106").Value = 21106
$ws.Range("L106").Value = 21106
$ws.Range("N106").Value = -23630
$ws.Range("H116").Value = 6999.9
$ws.Range("I116").Value = 6249.875
$ws.Range("K116").Value = 6249.875
$ws.Range("M116").Value = -3955.875
$ws.Range("H132").Value = 4918.0293
$ws.Range("I132").Value = 4823.6523
$ws.Range("J132").Value = 5115.364
$ws.Range("K132").Value = 14470.9569
$ws.Range("L132").Value = 15346.092
$ws.Range("M132").Value = -11940.9569
$ws.Range("N132").Value = -20406.092

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6999.9
$ws.Range("I3").Value = 6249.875
$ws.Range("K3").Value = 6249.875
$ws.Range("M3").Value = -6135.875
$ws.Range("H64").Value = 537.1
$ws.Range("J64").Value = 573.125
$ws.Range("L64").Value = 573.125
$ws.Range("N64").Value = -1023.125
$ws.Range("H67").Value = 537.1
$ws.Range("J67").Value = 573.125
$ws.Range("L67").Value = 573.125
$ws.Range("N67").Value = -2133.125
$ws.Range("H86").Value = 133333660
$ws.Range("H89").Value = 133333660
$ws.Range("H99").Value = 5486.25
$ws.Range("I99").Value = 3400
$ws.Range("J99").Value = 6529.375
$ws.Range("K99").Value = 3400
$ws.Range("L99").Value = 6529.375
$ws.Range("M99").Value = -1902
$ws.Range("N99").Value = -9525.375
$ws.Range("H105").Value = 3901.7222
$ws.Range("I105").Value = 3731.5715
$ws.Range("K105").Value = 3731.5715
$ws.Range("M105").Value = -1984.5715
$ws.Range("H134").Value = 5104.9
$ws.Range("I134").Value = 5033.222
$ws.Range("K134").Value = 15099.666
$ws.Range("M134").Value = -12564.666

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 9356.429
$ws.Range("I58").Value = 7110
$ws.Range("K58").Value = 7110
$ws.Range("M58").Value = -6907
$ws.Range("H86").Value = 5596.9546
$ws.Range("I86").Value = 3307.8235
$ws.Range("K86").Value = 3307.8235
$ws.Range("M86").Value = -2184.8235
$ws.Range("H89").Value = 5596.9546
$ws.Range("I89").Value = 3307.8235
$ws.Range("K89").Value = 16539.1175
$ws.Range("M89").Value = -10923.1175
$ws.Range("H94").Value = 1045.4117
$ws.Range("I94").Value = 777.4
$ws.Range("J94").Value = 1157.0834
$ws.Range("K94").Value = 777.4
$ws.Range("L94").Value = 1157.0834
$ws.Range("M94").Value = -326.4
$ws.Range("N94").Value = -2059.0834
$ws.Range("H132").Value = 6273.5454
$ws.Range("I132").Value = 5156.7144
$ws.Range("K132").Value = 15470.1432
$ws.Range("M132").Value = -12940.1432
$ws.Range("H134").Value = 7244.96
$ws.Range("I134").Value = 2844.0833
$ws.Range("K134").Value = 8532.249899999999
$ws.Range("M134").Value = -5997.249899999999
$ws.Range("H136").Value = 9356.429
$ws.Range("I136").Value = 7110
$ws.Range("K136").Value = 21330
$ws.Range("M136").Value = -18780

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 862.25
$ws.Range("I8").Value = 862.25
$ws.Range("K8").Value = 2586.75
$ws.Range("M8").Value = -2447.75
$ws.Range("H23").Value = 9394.75
$ws.Range("J23").Value = 11208.7
$ws.Range("L23").Value = 33626.10000000001
$ws.Range("N23").Value = -34096.10000000001
$ws.Range("H103").Value = 556
$ws.Range("I103").Value = 625
$ws.Range("J103").Value = 487
$ws.Range("K103").Value = 1875
$ws.Range("L103").Value = 1461
$ws.Range("M103").Value = -996
$ws.Range("N103").Value = -3219
$ws.Range("H140").Value = 2458.4375
$ws.Range("I140").Value = 1831.1666
$ws.Range("K140").Value = 5493.4998
$ws.Range("M140").Value = -313.4997999999996

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 785.1923
$ws.Range("I97").Value = 817.6
$ws.Range("K97").Value = 817.6
$ws.Range("M97").Value = -321.6
$ws.Range("H126").Value = 2925.25
$ws.Range("I126").Value = 3030.889
$ws.Range("K126").Value = 9092.667000000001
$ws.Range("M126").Value = -6622.667000000001
$ws.Range("H132").Value = 4871.3184
$ws.Range("I132").Value = 4624.737
$ws.Range("K132").Value = 13874.211
$ws.Range("M132").Value = -11344.211

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 14499.5
$ws.Range("I23").Value = 14499.5
$ws.Range("K23").Value = 14499.5
$ws.Range("M23").Value = -14269.5
$ws.Range("H103").Value = 19109.092
$ws.Range("J103").Value = 19109.092
$ws.Range("L103").Value = 19109.092
$ws.Range("N103").Value = -21453.092
$ws.Range("H122").Value = 5140.727
$ws.Range("J122").Value = 7358.1665
$ws.Range("L122").Value = 22074.4995
$ws.Range("N122").Value = -26974.4995
$ws.Range("H132").Value = 7891.68
$ws.Range("I132").Value = 9597
$ws.Range("K132").Value = 28791
$ws.Range("M132").Value = -26261
$ws.Range("H136").Value = 2722.7693
$ws.Range("I136").Value = 2656.8572
$ws.Range("J136").Value = 2799.6667
$ws.Range("K136").Value = 7970.571599999999
$ws.Range("L136").Value = 8399.000100000001
$ws.Range("M136").Value = -5420.571599999999
$ws.Range("N136").Value = -13499.0001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H41").Value = 16015
$ws.Range("J41").Value = 16015
$ws.Range("L41").Value = 16015
$ws.Range("N41").Value = -16795
$ws.Range("H113").Value = 760
$ws.Range("I113").Value = 717.1539
$ws.Range("K113").Value = 2151.4617
$ws.Range("M113").Value = 18.53830000000016
$ws.Range("H123").Value = 149499
$ws.Range("I123").Value = 149499
$ws.Range("K123").Value = 149499
$ws.Range("M123").Value = -144599
$ws.Range("H126").Value = 6307.769
$ws.Range("I126").Value = 3818.6365
$ws.Range("J126").Value = 19998
$ws.Range("K126").Value = 11455.9095
$ws.Range("L126").Value = 59994
$ws.Range("M126").Value = -8985.9095
$ws.Range("N126").Value = -64934
$ws.Range("H132").Value = 4429.8335
$ws.Range("I132").Value = 3380.6155
$ws.Range("J132").Value = 11249.75
$ws.Range("K132").Value = 10141.8465
$ws.Range("L132").Value = 33749.25
$ws.Range("M132").Value = -7611.8465
$ws.Range("N132").Value = -38809.25
